# Apply updated cryptocurrency price/volume data (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.347.81"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.936.05"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'252.05"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("D6").Value = "'0.7247"
$ws.Range("E6").Value = "  +4.07%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.3310"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").Value = "'27.92"
$ws.Range("E9").Value = "  +6.86%  "
$ws.Range("D10").Value = "'0.07249"
$ws.Range("E10").Value = "  +6.71%  "
$ws.Range("D11").Value = "'0.8092"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "'0.08105"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Value = "1.935.18"
$ws.Range("D14").Value = "'5.481"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "'94.88"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'15.11"
$ws.Range("E16").Value = "  +5.22%  "
$ws.Range("D17").Value = "30.344.34"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "'0.000008221"
$ws.Range("E18").Value = "  +5.38%  "
$ws.Range("D19").Value = "'253.03"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").Value = "'5.836"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "2.190.62"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "'6.962"
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("D25").Value = "'9.761"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").Value = "'166.07"
$ws.Range("E26").Value = "  +3.95%  "
$ws.Range("E27").Value = "  +6.52%  "
$ws.Range("D28").Value = "'19.34"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("D29").Value = "'0.1301"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'1.546"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'4.440"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").Value = "'4.217"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").Value = "'0.05249"
$ws.Range("E34").Value = "  +4.51%  "
$ws.Range("E35").Value = "  +7.21%  "
$ws.Range("D36").Value = "'0.7513"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("D37").Value = "'2.774"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").Value = "'0.01971"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").Value = "'2.806"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").Value = "'79.43"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'6.452"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "'0.4550"
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").Value = "'0.8440"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'101.99"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'9.829"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").Value = "'7.448"
$ws.Range("E48").Value = "  +4.25%  "
$ws.Range("D49").Value = "'36.80"
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("D50").Value = "'0.4202"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("D51").Value = "'0.06050"
$ws.Range("E51").Value = "  +2.18%  "
